$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "선형 연립 미분방정식 모델링"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/11/modeling_with_systems.html"

$ws.Range("D8").Value = "제주어 기계번역 모델과 음성합성 모델에 관한 연구를 소개합니다."

$ws.Range("D28").Value = "[강화학습] Markov Decision Processes(2)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/118"

$ws.Range("D29").Value = "[만화] 인턴일기 19~27"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-4/"

$ws.Range("D37").Value = "[Paper Review] SSD: A unified framework for self-supervised outlier detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1463&mod=document&pageid=1"

$ws.Range("D50").Value = "인공지능 기술로 어디까지 넘볼 수 있나?"
$ws.Range("E50").Value = "http://incredible.egloos.com/7506674"

$ws.Range("D51").Value = "[html] 5초마다 자동으로 웹 페이지에 새로고침을 해주고 싶다면?"
$ws.Range("E51").Value = "https://bskyvision.com/1192"
